# Atualiza os resultados da tabela de medicamentos com os novos dados.
$d = $word.ActiveDocument

# n
$d.Content.Find.Execute("422", $true, $false, $false, $false, $false, $true, 1, $false, "365", 2)

# Profilaxia (%)
$d.Content.Find.Execute("388 (95.1)", $true, $false, $false, $false, $false, $true, 1, $false, "321 (96,4)", 2)

# Dabigatrana (%)
$d.Content.Find.Execute("5 ( 1.2)", $true, $false, $false, $false, $false, $true, 1, $false, "5 ( 1,5)", 2)

# Enoxaparina (%)
$d.Content.Find.Execute("380 (93.1)", $true, $false, $false, $false, $false, $true, 1, $false, "319 (95,5)", 2)

# Corrige o nome "Rivoraxabana" para "Rivaroxabana"
$d.Content.Find.Execute("Rivoraxabana (%)", $true, $false, $false, $false, $false, $true, 1, $false, "Rivaroxabana (%)", 2)

# Rivaroxabana (%) - quantidade
$d.Content.Find.Execute("86 (20.9)", $true, $false, $false, $false, $false, $true, 1, $false, "72 (21,5)", 2)

# Warfarina (%)
$d.Content.Find.Execute("193 (51.7)", $true, $false, $false, $false, $false, $true, 1, $false, "177 (57,5)", 2)
